$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.423.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.884.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("E4").Value = "  -0.77%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.74%  "
$ws.Range("E6").Value = "  -4.90%  "
$ws.Range("E7").Value = "  -0.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.01"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.03%  "
$ws.Range("E9").Value = "  -4.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0736"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0969"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "13.04"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.158.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.738"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("E15").Value = "  -1.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.871.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.382.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "73.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0820"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "244.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.50%  "
$ws.Range("E22").Value = "  -4.59%  "
$ws.Range("E23").Value = "  -0.82%  "
$ws.Range("E24").Value = "  +3.05%  "
$ws.Range("E25").Value = "  -9.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("E27").Value = "  -3.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.56%  "
$ws.Range("E29").Value = "  -4.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.128.43"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  +2.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0577"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.51%  "
$ws.Range("E34").Value = "  -2.13%  "
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -13.26%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.851"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.45%  "
$ws.Range("E38").Value = "  -4.26%  "
$ws.Range("E39").Value = "  +4.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "97.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.03%  "
$ws.Range("E41").Value = "  -2.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("E43").Value = "  -4.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.291.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.01%  "
$ws.Range("E45").Value = "  -5.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0809"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.65%  "
$ws.Range("E47").Value = "  -1.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "43.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.26%  "
$ws.Range("E51").Value = "  -7.63%  "
